$wb = $excel.ActiveWorkbook

# --- Hangar sheet: move selection from G8 to B9 ---
$hangar = $wb.Worksheets.Item("Hangar")
$hangar.Activate()
$hangar.Range("B9").Select()

# --- Empleado sheet: fill in the real menu options (placeholder text was "Empleado") ---
# Write order matters for shared-string table ordering, so write B3..B6 first,
# then B2, then B7, matching the target shared-strings layout.
$empleado = $wb.Worksheets.Item("Empleado")
$empleado.Range("B3").Value = "Nuevo empleado"
$empleado.Range("B4").Value = "Actualizar informacion del empleado"
$empleado.Range("B5").Value = "Actualizar el estado del empleado"
$empleado.Range("B6").Value = "Verificar disponibilidad para el vuelo"
$empleado.Range("B2").Value = "Informacion del empleado"
$empleado.Range("B7").Value = "Regresar al menu principal"

# Make Empleado the active sheet/tab with B2 selected.
$empleado.Activate()
$empleado.Range("B2").Select()
